$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C, rows 2 through 33 hold a date (style "YYYY-MM-DD") that was
# bumped forward by one day (45189 -> 45190, i.e. 2023-09-20 -> 2023-09-21).
for ($row = 2; $row -le 33; $row++) {
    $ws.Cells.Item($row, 3).Value = 45190
}
